# Identity and report fix
# Update "Количество категорий" (C), "Количество скриншотов программы" (D)
# and "Количество комментариев программы" (E) counts for each software row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = @(3, 4, 2)
    4  = @(1, 3, 4)
    5  = @(1, 0, 4)
    6  = @(1, 0, 4)
    7  = @(1, 0, 4)
    8  = @(1, 0, 4)
    9  = @(1, 0, 4)
    10 = @(1, 0, 4)
    11 = @(1, 0, 4)
    12 = @(1, 0, 4)
    13 = @(1, 0, 4)
    14 = @(1, 0, 4)
    15 = @(1, 0, 4)
    16 = @(1, 0, 4)
    17 = @(1, 0, 4)
    18 = @(1, 0, 4)
    19 = @(1, 0, 4)
    20 = @(1, 0, 4)
}

foreach ($row in $values.Keys | Sort-Object) {
    $vals = $values[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
}
